$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3401.8386
$ws.Range("I33").Value = 4843.4287
$ws.Range("J33").Value = 374.5
$ws.Range("K33").Value = 4843.4287
$ws.Range("L33").Value = 374.5
$ws.Range("M33").Value = -4614.4287
$ws.Range("N33").Value = -832.5
$ws.Range("H64").Value = 6514.2856
$ws.Range("I64").Value = 7860
$ws.Range("J64").Value = 3150
$ws.Range("K64").Value = 7860
$ws.Range("L64").Value = 3150
$ws.Range("M64").Value = -7612
$ws.Range("N64").Value = -3646
$ws.Range("H67").Value = 6514.2856
$ws.Range("I67").Value = 7860
$ws.Range("J67").Value = 3150
$ws.Range("K67").Value = 7860
$ws.Range("L67").Value = 3150
$ws.Range("M67").Value = -7002
$ws.Range("N67").Value = -4866
$ws.Range("H112").Value = 3622.2222
$ws.Range("J112").Value = 3723.8096
$ws.Range("L112").Value = 11171.4288
$ws.Range("N112").Value = -13387.4288
$ws.Range("H118").Value = 772.4167
$ws.Range("I118").Value = 524.4545000000001
$ws.Range("K118").Value = 1573.3635
$ws.Range("M118").Value = 83.63649999999984
$ws.Range("H121").Value = 446.94
$ws.Range("J121").Value = 432.22916
$ws.Range("L121").Value = 1296.68748
$ws.Range("N121").Value = -4790.687480000001
$ws.Range("H129").Value = 2179797.2
$ws.Range("J129").Value = 3705343.5
$ws.Range("L129").Value = 11116030.5
$ws.Range("N129").Value = -11126030.5
$ws.Range("H132").Value = 2042666.2
$ws.Range("I132").Value = 1777.7632
$ws.Range("J132").Value = 9093008
$ws.Range("K132").Value = 5333.2896
$ws.Range("L132").Value = 27279024
$ws.Range("M132").Value = -2803.2896
$ws.Range("N132").Value = -27284084
$ws.Range("H137").Value = 1014.24
$ws.Range("I137").Value = 886.5227
$ws.Range("J137").Value = 1950.8334
$ws.Range("K137").Value = 2659.5681
$ws.Range("L137").Value = 5852.5002
$ws.Range("M137").Value = -109.5681
$ws.Range("N137").Value = -10952.5002
$ws.Range("H138").Value = 3835.7024
$ws.Range("I138").Value = 1850.5667
$ws.Range("J138").Value = 4938.5557
$ws.Range("K138").Value = 5551.7001
$ws.Range("L138").Value = 14815.6671
$ws.Range("M138").Value = -411.7001
$ws.Range("N138").Value = -25095.6671

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18995.451
$ws.Range("I32").Value = 14903.728
$ws.Range("J32").Value = 82008
$ws.Range("K32").Value = 14903.728
$ws.Range("L32").Value = 82008
$ws.Range("M32").Value = -14616.728
$ws.Range("N32").Value = -82582
$ws.Range("H45").Value = 1326.5
$ws.Range("I45").Value = 1383
$ws.Range("J45").Value = 1251.1666
$ws.Range("K45").Value = 1383
$ws.Range("L45").Value = 1251.1666
$ws.Range("M45").Value = -1006
$ws.Range("N45").Value = -2005.1666
$ws.Range("H110").Value = 741.94446
$ws.Range("I110").Value = 703.73334
$ws.Range("J110").Value = 933
$ws.Range("K110").Value = 703.73334
$ws.Range("L110").Value = 933
$ws.Range("M110").Value = 1341.26666
$ws.Range("N110").Value = -5023

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5808.5
$ws.Range("I107").Value = 719.13635
$ws.Range("J107").Value = 33800
$ws.Range("K107").Value = 719.13635
$ws.Range("L107").Value = 33800
$ws.Range("M107").Value = 1200.86365
$ws.Range("N107").Value = -37640
$ws.Range("H134").Value = 20789.434
$ws.Range("I134").Value = 1861.4651
$ws.Range("J134").Value = 102179.7
$ws.Range("K134").Value = 5584.3953
$ws.Range("L134").Value = 306539.1
$ws.Range("M134").Value = -3049.3953
$ws.Range("N134").Value = -311609.1

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2108.6924
$ws.Range("I16").Value = 1083.3334
$ws.Range("J16").Value = 2987.5715
$ws.Range("K16").Value = 1083.3334
$ws.Range("L16").Value = 2987.5715
$ws.Range("M16").Value = -796.3334
$ws.Range("N16").Value = -3561.5715
$ws.Range("H31").Value = 2987.9517
$ws.Range("I31").Value = 2345.795
$ws.Range("J31").Value = 4076.8262
$ws.Range("K31").Value = 2345.795
$ws.Range("L31").Value = 4076.8262
$ws.Range("M31").Value = -2050.795
$ws.Range("N31").Value = -4666.8262
$ws.Range("H34").Value = 2987.9517
$ws.Range("I34").Value = 2345.795
$ws.Range("J34").Value = 4076.8262
$ws.Range("K34").Value = 2345.795
$ws.Range("L34").Value = 4076.8262
$ws.Range("M34").Value = -2143.795
$ws.Range("N34").Value = -4480.8262
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H113").Value = 2108.6924
$ws.Range("I113").Value = 1083.3334
$ws.Range("J113").Value = 2987.5715
$ws.Range("K113").Value = 1083.3334
$ws.Range("L113").Value = 2987.5715
$ws.Range("M113").Value = 1086.6666
$ws.Range("N113").Value = -7327.5715

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2172.9565
$ws.Range("I5").Value = 1730.6666
$ws.Range("J5").Value = 2329.0588
$ws.Range("K5").Value = 5191.9998
$ws.Range("L5").Value = 6987.176399999999
$ws.Range("M5").Value = -5079.9998
$ws.Range("N5").Value = -7211.176399999999
$ws.Range("H23").Value = 124.6875
$ws.Range("I23").Value = 66
$ws.Range("J23").Value = 159.9
$ws.Range("K23").Value = 198
$ws.Range("L23").Value = 479.7
$ws.Range("M23").Value = 37
$ws.Range("N23").Value = -949.7
$ws.Range("H115").Value = 1188.2858
$ws.Range("I115").Value = 509.33334
$ws.Range("J115").Value = 1697.5
$ws.Range("K115").Value = 1528.00002
$ws.Range("L115").Value = 5092.5
$ws.Range("M115").Value = -353.0000199999999
$ws.Range("N115").Value = -7442.5
$ws.Range("H122").Value = 1448.4286
$ws.Range("J122").Value = 1631
$ws.Range("L122").Value = 14679
$ws.Range("N122").Value = -19579
$ws.Range("H123").Value = 1764.6666
$ws.Range("I123").Value = 1033.3334
$ws.Range("J123").Value = 2496
$ws.Range("K123").Value = 3100.0002
$ws.Range("L123").Value = 7488
$ws.Range("M123").Value = -650.0001999999999
$ws.Range("N123").Value = -12388
$ws.Range("H135").Value = 2172.9565
$ws.Range("I135").Value = 1730.6666
$ws.Range("J135").Value = 2329.0588
$ws.Range("K135").Value = 15575.9994
$ws.Range("L135").Value = 20961.5292
$ws.Range("M135").Value = -13040.9994
$ws.Range("N135").Value = -26031.5292
$ws.Range("H137").Value = 52948.15
$ws.Range("I137").Value = 3467.2727
$ws.Range("J137").Value = 113424.78
$ws.Range("K137").Value = 10401.8181
$ws.Range("L137").Value = 340274.34
$ws.Range("M137").Value = -5301.8181
$ws.Range("N137").Value = -350474.34

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H126").Value = 2632.6
$ws.Range("I126").Value = 1972.238
$ws.Range("J126").Value = 6099.5
$ws.Range("K126").Value = 5916.714
$ws.Range("L126").Value = 18298.5
$ws.Range("M126").Value = -3446.714
$ws.Range("N126").Value = -23238.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").ClearContents()
$ws.Range("H82").Value = 2846.4707
$ws.Range("I82").Value = 2478.1428
$ws.Range("J82").Value = 4565.3335
$ws.Range("K82").Value = 2478.1428
$ws.Range("L82").Value = 4565.3335
$ws.Range("M82").Value = -2117.1428
$ws.Range("N82").Value = -5287.3335
$ws.Range("H85").Value = 2846.4707
$ws.Range("I85").Value = 2478.1428
$ws.Range("J85").Value = 4565.3335
$ws.Range("K85").Value = 2478.1428
$ws.Range("L85").Value = 4565.3335
$ws.Range("M85").Value = -1230.1428
$ws.Range("N85").Value = -7061.3335

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 10000
$ws.Range("K54").Value = 10000
$ws.Range("M54").Value = -9480
$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2692
